$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.966.26'
$ws.Range('E2').Value = '  +5.07%  '
$ws.Range('D3').Value = '3.516.61'
$ws.Range('E3').Value = '  +2.75%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '592.74'
$ws.Range('E5').Value = '  +3.94%  '
$ws.Range('D6').Value = '168.83'
$ws.Range('E6').Value = '  +6.69%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.513.70'
$ws.Range('E8').Value = '  +2.52%  '
$ws.Range('D9').Value = '0.578'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').Value = '7.29'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('E11').Value = '  +5.48%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.440'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.07%  '
$ws.Range('D13').Value = '4.120.66'
$ws.Range('E13').Value = '  +2.69%  '
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('D15').Value = '28.21'
$ws.Range('E15').Value = '  +4.13%  '
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').Value = '66.932.82'
$ws.Range('E17').Value = '  +4.86%  '
$ws.Range('D18').Value = '3.520.82'
$ws.Range('E18').Value = '  +3.06%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  +4.14%  '
$ws.Range('D20').Value = '14.02'
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('D21').Value = '391.52'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').Value = '73.65'
$ws.Range('E23').Value = '  +3.53%  '
$ws.Range('E24').Value = '  +10.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = '0.533'
$ws.Range('E26').Value = '  +3.56%  '
$ws.Range('D27').Value = '10.19'
$ws.Range('E27').Value = '  +5.63%  '
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('E29').Value = '  +0.88%  '
$ws.Range('D30').Value = '6.43'
$ws.Range('E30').Value = '  +5.86%  '
$ws.Range('D31').Value = '1.48'
$ws.Range('E31').Value = '  +6.44%  '
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  +3.98%  '
$ws.Range('D33').Value = '23.61'
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('D34').Value = '7.48'
$ws.Range('E34').Value = '  +7.59%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E36').Value = '  +5.94%  '
$ws.Range('D37').Value = '161.03'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  +6.09%  '
$ws.Range('E39').Value = '  +5.55%  '
$ws.Range('E40').Value = '  +3.74%  '
$ws.Range('E41').Value = '  +7.35%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '6.71'
$ws.Range('E42').Value = '  +5.04%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '26.49'
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('D44').Value = '2.832.83'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '43.52'
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('D46').Value = '26.41'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0315'
$ws.Range('E47').Value = '  +3.73%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').Value = '2.54'
$ws.Range('E48').Value = '  +5.74%  '
$ws.Range('D49').Value = '353.28'
$ws.Range('E49').Value = '  +6.14%  '
$ws.Range('E50').Value = '  +4.45%  '
$ws.Range('D51').Value = '33.63'
$ws.Range('E51').Value = '  +12.40%  '
